# Rewrites the first two npm-related bullet points in the document:
#   1) "npm init - makes ..."      -> quotes `npm init`, "save and update" -> "save all"
#   2) "npm install ... - Version is optional" -> quotes the full command (with --save),
#      and appends an explanation of the --save flag / package.json auto-update.
$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$target1 = $d.Range($r1.Start, $r1.End - 1)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>“</w:t></w:r><w:r><w:t>npm init</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space="preserve"> – makes the package.json</w:t></w:r><w:r><w:t xml:space="preserve"> file that will save all </w:t></w:r><w:r><w:t>dependencies</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target1.InsertXML($xml1)

$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$target2 = $d.Range($r2.Start, $r2.End - 1)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>“</w:t></w:r><w:r><w:t xml:space="preserve">npm install &lt;package name&gt;@&lt;version&gt; </w:t></w:r><w:r><w:t xml:space="preserve">--save” </w:t></w:r><w:r><w:t>- Version is optional</w:t></w:r><w:r><w:t xml:space="preserve"> and –save flag saves this new dependency/package name i</w:t></w:r><w:r><w:t>n pack</w:t></w:r><w:r><w:t>age.json</w:t></w:r><w:r><w:t>. New version of node probably have automatic updation for package.json file, so no need for –save flag.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target2.InsertXML($xml2)

Write-Host "Paragraph 1 now: " $d.Paragraphs(1).Range.Text
Write-Host "Paragraph 2 now: " $d.Paragraphs(2).Range.Text
